{"js": "// Fix typo: \" Baixo \" -> \" Baixa \" (gender agreement with \"Estamina\")\n// in the \"Estamina e Moral\" paragraph of the SDM Quick Guide.\nconst body = context.document.body;\n\n// Match the exact original run text \" Baixo \" (leading + trailing space),\n// case-sensitive, so we don't touch the unrelated lowercase \"baixo\" earlier\n// in the same paragraph (\"quando o moral estiver baixo\").\nconst results = body.search(\" Baixo \", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the text ' Baixo ' to correct.\");\n}\n\n// Replace in place, preserving the surrounding single leading/trailing space.\nresults.items[0].insertText(\" Baixa \", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix typo: \" Baixo \" -> \" Baixa \" (gender agreement with \"Estamina\")\n# in the \"Estamina e Moral\" paragraph of the SDM Quick Guide.\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n# Case-sensitive match on \" Baixo \" (leading + trailing space) so the\n# unrelated lowercase \"baixo\" earlier in the same paragraph (\"quando o\n# moral estiver baixo\") is left untouched.\n$found = $rng.Find.Execute(\n    \" Baixo \",   # FindText\n    $true,       # MatchCase\n    $false,      # MatchWholeWord\n    $false,      # MatchWildcards\n    $false,      # MatchSoundsLike\n    $false,      # MatchAllWordForms\n    $true,       # Forward\n    1,           # Wrap (wdFindContinue)\n    $false,      # Format\n    \" Baixa \",   # ReplaceWith\n    2            # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw \"Could not find the text ' Baixo ' to correct.\"\n}\n"}
